$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the tiny floating-point discrepancy already present in B53
$ws.Cells.Item(53, 2).Value = 45735.96277708333

# The date/time number format already used by the existing B column
$dateFmt = $ws.Cells.Item(53, 2).NumberFormat

# New "Marvin" login timestamps appended as rows 54-75
$newRows = @(
    45735.96462568287,
    45735.96465396991,
    45735.96465576389,
    45735.96465761574,
    45735.96465920139,
    45735.96466092593,
    45735.96466273148,
    45735.96466434028,
    45735.96466599537,
    45735.96466791667,
    45735.96467282408,
    45735.96467512732,
    45735.96467731481,
    45735.96468119213,
    45735.96468372685,
    45735.96468569445,
    45735.96470547454,
    45735.9647071875,
    45735.96473335648,
    45735.96475530093,
    45735.96477335648,
    45735.96752236925
)

$row = 54
foreach ($ts in $newRows) {
    $ws.Cells.Item($row, 1).Value = "Marvin"
    $ws.Cells.Item($row, 2).Value = $ts
    $ws.Cells.Item($row, 2).NumberFormat = $dateFmt
    $row++
}
